# Auto-update price data: insert a new "latest day" row at the top of the
# data table (row 2), pushing all existing rows down by one. The new row
# carries the same commodity values as the previous top row and a date that
# is one day later than the previous most-recent date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (the first data row, below the header),
# shifting everything below it (and the old row 54 at the bottom) down by
# one row.
$ws.Rows.Item(2).Insert()

# A fresh row-insert can pick up neighbouring formatting (e.g. the bold
# header style); strip that so the new row matches the plain, unstyled
# data rows around it.
$ws.Rows.Item(2).ClearFormats()

# Force the date cell to be stored as plain text (matching every other
# date cell in the column) instead of letting Excel auto-convert the
# "yyyy-mm-dd" literal into a date serial number; ClearFormats()
# afterwards drops the temporary text number-format again so the cell
# keeps the sheet's default (unstyled) look.
$dateCell = $ws.Cells.Item(2, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-01-13"
$dateCell.ClearFormats()

# Same commodity values as every other row in this (currently flat) series.
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
